# issue #5: add legislator_id, name, date into dataframe
# The "股票" (stocks) sheet gets three new trailing columns:
#   H = date            (text, e.g. "2012-04-20")
#   I = legislator_name (e.g. "費鴻泰")
#   J = legislator_id   (numeric id, e.g. 1365)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Force column H to be treated as text so the date string is not
# auto-converted into a date serial number.
$ws.Range("H1:H7").NumberFormat = "@"

# Header row (row 1)
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Data rows (rows 2-7)
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 8).Value = "2012-04-20"
    $ws.Cells.Item($r, 9).Value = "費鴻泰"
    $ws.Cells.Item($r, 10).Value = 1365
}

# Copy the formatting of the existing header/data cells onto the new
# columns so they pick up the same styling (bold+border header, plain
# data rows) as the rest of the table.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 7).Copy()
    $dst = $ws.Range($ws.Cells.Item($r, 8), $ws.Cells.Item($r, 10))
    $dst.PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
